$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so that numeric-looking
# strings (e.g. "5.41", "23.42") are not auto-converted into real numbers
# by the automatic type inference, preserving the original text values
# and formatting exactly as they appear in the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '57.877.82'
$ws.Range('E2').Value = '  -3.53%  '
$ws.Range('D3').Value = '2.284.64'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '532.66'
$ws.Range('E5').Value = '  -4.45%  '
$ws.Range('D6').Value = '130.42'
$ws.Range('E6').Value = '  -2.71%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.581'
$ws.Range('E8').Value = '  -0.88%  '
$ws.Range('D9').Value = '2.282.06'
$ws.Range('E9').Value = '  -4.15%  '
$ws.Range('D10').Value = '0.0993'
$ws.Range('E10').Value = '  -6.51%  '
$ws.Range('D11').Value = '5.41'
$ws.Range('E11').Value = '  -4.48%  '
$ws.Range('E12').Value = '  -0.55%  '
$ws.Range('E13').Value = '  -4.30%  '
$ws.Range('D14').Value = '23.42'
$ws.Range('E14').Value = '  -4.31%  '
$ws.Range('D15').Value = '2.693.15'
$ws.Range('E15').Value = '  -4.04%  '
$ws.Range('D16').Value = '57.817.47'
$ws.Range('E16').Value = '  -3.61%  '
$ws.Range('D17').Value = '0.0000130'
$ws.Range('E17').Value = '  -5.42%  '
$ws.Range('D18').Value = '2.282.18'
$ws.Range('E18').Value = '  -3.91%  '
$ws.Range('D19').Value = '10.46'
$ws.Range('E19').Value = '  -6.05%  '
$ws.Range('D20').Value = '4.21'
$ws.Range('E20').Value = '  -6.48%  '
$ws.Range('D21').Value = '311.20'
$ws.Range('E21').Value = '  -3.24%  '
$ws.Range('E22').Value = '  -4.81%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').Value = '62.27'
$ws.Range('E24').Value = '  -2.91%  '
$ws.Range('E25').Value = '  -3.98%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').Value = '7.97'
$ws.Range('E27').Value = '  -5.66%  '
$ws.Range('E28').Value = '  -7.24%  '
$ws.Range('D29').Value = '170.44'
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('E30').Value = '  -6.20%  '
$ws.Range('D31').Value = '0.0₃0715'
$ws.Range('E31').Value = '  -6.13%  '
$ws.Range('D32').Value = '5.73'
$ws.Range('E32').Value = '  -5.94%  '
$ws.Range('E33').Value = '  -7.32%  '
$ws.Range('E34').Value = '  -5.78%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').Value = '17.67'
$ws.Range('E36').Value = '  -2.65%  '
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('E38').Value = '  -7.90%  '
$ws.Range('D39').Value = '3.87'
$ws.Range('E39').Value = '  -6.78%  '
$ws.Range('D40').Value = '38.12'
$ws.Range('E40').Value = '  -1.33%  '
$ws.Range('E41').Value = '  -6.95%  '
$ws.Range('D42').Value = '141.18'
$ws.Range('E42').Value = '  -2.96%  '
$ws.Range('D43').Value = '286.08'
$ws.Range('E43').Value = '  -10.43%  '
$ws.Range('E44').Value = '  -3.89%  '
$ws.Range('D45').Value = '0.0943'
$ws.Range('E45').Value = '  -2.98%  '
$ws.Range('D46').Value = '0.0493'
$ws.Range('E46').Value = '  -3.82%  '
$ws.Range('D47').Value = '0.552'
$ws.Range('E47').Value = '  -3.17%  '
$ws.Range('D48').Value = '18.03'
$ws.Range('E48').Value = '  -8.94%  '
$ws.Range('D49').Value = '0.0209'
$ws.Range('E49').Value = '  -4.05%  '
$ws.Range('E50').Value = '  -1.10%  '
$ws.Range('B51').Value = 'ZEEBU'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range('D51').Value = '4.64'
$ws.Range('E51').Value = '  -0.72%  '
